# Financials update: insert a new first data column (FY ending 2018-12-31,
# serial date 43465) in front of the existing "Period Ending" columns on the
# FXNC sheet, pushing the prior years one column to the right, and populate
# the new column with the latest figures for every line item across the
# Income Statement, Balance Sheet and Cash Flow Statement blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D; Excel shifts D:K -> E:L automatically (formulas,
# dimension and row "spans" all update as part of the native insert).
$ws.Columns("D:D").Insert()

# The freshly inserted column D inherits formatting from the column to its
# left (C, the label column). Re-apply the per-row number format/font that
# the rest of the data columns use by copying formats from column E (which
# now holds what used to be column D, so it already carries the correct
# per-row style - date format on the "Period Ending" rows, number format
# elsewhere).
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 37 and 79 are bare section headers (only column B holds text) with no
# data columns at all; PasteSpecial stamped an empty, styled D cell onto them
# even though neither the row before nor after the edit has one there, so
# drop it again.
$ws.Range("D37").Clear() | Out-Null
$ws.Range("D79").Clear() | Out-Null

# --- Income Statement -------------------------------------------------
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 31100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -1800
$ws.Range("D17").Value = 4100
$ws.Range("D18").Value = 27000
$ws.Range("D20").Value = -14600
$ws.Range("D21").Value = 14200
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 12400
$ws.Range("D24").Value = 2300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 10100
$ws.Range("D27").Value = 10100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 14600
$ws.Range("D33").Value = 10100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 10100

# --- Balance Sheet ------------------------------------------------------
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 13400
$ws.Range("D42").Value = 16900
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 20100
$ws.Range("D49").Value = 500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 753000
$ws.Range("D57").Value = 1500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 14200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 686300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 54800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 66700
$ws.Range("D77").Value = 0

# --- Cash Flow Statement -------------------------------------------------
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 10100
$ws.Range("D83").Value = 1800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 13800
$ws.Range("D91").Value = -1500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -29800
$ws.Range("D96").Value = -900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 4700
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -11400
